$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card20")

# Copy formatting (style/borders) from column M down to the new column N
# for the full used range (header + data rows) BEFORE changing any values,
# so the new column N picks up the same per-row formatting as column M
# (bold/bordered header in row 1, plain/unstyled cells in rows 2-12).
$ws.Range("M1:M12").Copy()
$ws.Range("N1:N12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header: "Event " -> "Event" (drop trailing space)
$ws.Cells.Item(1, 13).Value = "Event"

# New header for column N
$ws.Cells.Item(1, 14).Value = "Correction "

# Column M data rows: blank cells become the text "nan"
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 13).Value = "nan"
}

# Dimension / used range will now naturally extend to N12
